$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.1465433017255542
$ws.Range("D2").Value = 0.1262258172897148
$ws.Range("E2").Value = 0.1372242642186237
$ws.Range("F2").Value = 1.049929521006806
$ws.Range("G2").Value = 0.593079604595367
$ws.Range("H2").Value = 0.7286482229212794
$ws.Range("I2").Value = 0.7286526594982874
$ws.Range("J2").Value = 0.2105825628211022
$ws.Range("L2").Value = 0.2471504961990263
$ws.Range("M2").Value = 10.99935363126821
$ws.Range("O2").Value = 2.624169771381901
$ws.Range("C3").Value = 0.1512743028483996
$ws.Range("D3").Value = 0.1294733883543273
$ws.Range("E3").Value = 0.1358730885044181
$ws.Range("F3").Value = 1.081708089017781
$ws.Range("G3").Value = 0.6098825884774826
$ws.Range("H3").Value = 0.744692625321477
$ws.Range("I3").Value = 0.7506023444175973
$ws.Range("J3").Value = 0.2012473714441541
$ws.Range("L3").Value = 0.2363310134679182
$ws.Range("M3").Value = 9.672548105483997
$ws.Range("O3").Value = 2.692913045803508
$ws.Range("C4").Value = 0.1543393044362844
$ws.Range("D4").Value = 0.1315803123318169
$ws.Range("E4").Value = 0.1351913100944913
$ws.Range("F4").Value = 1.10252289669214
$ws.Range("G4").Value = 0.6211072812141225
$ws.Range("H4").Value = 0.7551910564542723
$ws.Range("I4").Value = 0.7649978706321505
$ws.Range("J4").Value = 0.1957148321517082
$ws.Range("L4").Value = 0.2298975661908145
$ws.Range("M4").Value = 8.854117297462267
$ws.Range("O4").Value = 2.738387094076174
$ws.Range("C5").Value = 0.1556282079298121
$ws.Range("D5").Value = 0.1324670792130256
$ws.Range("E5").Value = 0.1349500924337335
$ws.Range("F5").Value = 1.111329145802284
$ws.Range("G5").Value = 0.6259066778358573
$ws.Range("H5").Value = 0.7596307204188122
$ws.Range("I5").Value = 0.7710925688753889
$ws.Range("J5").Value = 0.1935093440413738
$ws.Range("L5").Value = 0.2273276267139721
$ws.Range("M5").Value = 8.519638314464089
$ws.Range("O5").Value = 2.757730306520671
$ws.Range("C6").Value = 0.1558446276446155
$ws.Range("D6").Value = 0.1326160219496657
$ws.Range("E6").Value = 0.1349122337891373
$ws.Range("F6").Value = 1.112810885721924
$ws.Range("G6").Value = 0.6267171311861119
$ws.Range("H6").Value = 0.760377638883817
$ws.Range("I6").Value = 0.7721183104135285
$ws.Range("J6").Value = 0.1931460585754934
$ws.Range("L6").Value = 0.2269039874870629
$ws.Range("M6").Value = 8.464039895500775
$ws.Range("O6").Value = 2.760991023950496
$ws.Range("C7").Value = 0.154356526022486
$ws.Range("D7").Value = 0.1315921577643984
$ws.Range("E7").Value = 0.1351879094298098
$ws.Range("F7").Value = 1.102640353424846
$ws.Range("G7").Value = 0.6211710994772872
$ws.Range("H7").Value = 0.7552502791865905
$ws.Range("I7").Value = 0.7650791443790474
$ws.Range("J7").Value = 0.1956848908350537
$ws.Range("L7").Value = 0.2298626988537507
$ws.Range("M7").Value = 8.849610307620992
$ws.Range("O7").Value = 2.738644686745928
$ws.Range("C8").Value = 0.1481410827165135
$ws.Range("D8").Value = 0.1273220187163204
$ws.Range("E8").Value = 0.136727342354952
$ws.Range("F8").Value = 1.060614297525248
$ws.Range("G8").Value = 0.5986831098797438
$ws.Range("H8").Value = 0.7340452763284091
$ws.Range("I8").Value = 0.7360287982724572
$ws.Range("J8").Value = 0.2073217653236128
$ws.Range("L8").Value = 0.2433757719745699
$ws.Range("M8").Value = 10.54264252438145
$ws.Range("O8").Value = 2.647189341890936
$ws.Range("C9").Value = 0.1372403487916785
$ws.Range("D9").Value = 0.1198536471905385
$ws.Range("E9").Value = 0.1409466561378352
$ws.Range("F9").Value = 0.9886933013169852
$ws.Range("G9").Value = 0.5619188535974189
$ws.Range("H9").Value = 0.6976520639155837
$ws.Range("I9").Value = 0.6864596482817333
$ws.Range("J9").Value = 0.2317753936245879
$ws.Range("L9").Value = 0.2715899386168275
$ws.Range("M9").Value = 13.83380611208622
$ws.Range("O9").Value = 2.494153857113247
$ws.Range("C10").Value = 0.1300400915553155
$ws.Range("D10").Value = 0.114931752576247
$ws.Range("E10").Value = 0.1448191452812821
$ws.Range("F10").Value = 0.9424622975217076
$ws.Range("G10").Value = 0.5395628461749595
$ws.Range("H10").Value = 0.6741527246731067
$ws.Range("I10").Value = 0.6547044702664202
$ws.Range("J10").Value = 0.2508172733684262
$ws.Range("L10").Value = 0.2934427726970199
$ws.Range("M10").Value = 16.23616699741677
$ws.Range("O10").Value = 2.398303039083146
$ws.Range("C11").Value = 0.1269452724409419
$ws.Range("D11").Value = 0.1128183102810709
$ws.Range("E11").Value = 0.1467581589733129
$ws.Range("F11").Value = 0.9229120762650069
$ws.Range("G11").Value = 0.530443158372492
$ws.Range("H11").Value = 0.6641817266102237
$ws.Range("I11").Value = 0.6413040351384467
$ws.Range("J11").Value = 0.2597325770297516
$ws.Range("L11").Value = 0.3036471109371917
$ws.Range("M11").Value = 17.32618580609409
$ws.Range("O11").Value = 2.358417832088008
$ws.Range("C12").Value = 0.1257997705094809
$ws.Range("D12").Value = 0.1120363180587916
$ws.Range("E12").Value = 0.1475187522745358
$ws.Range("F12").Value = 0.915725740366625
$ws.Range("G12").Value = 0.5271441430749064
$ws.Range("H12").Value = 0.6605107342893746
$ws.Range("I12").Value = 0.6363827498063195
$ws.Range("J12").Value = 0.2631466057509471
$ws.Range("L12").Value = 0.3075507326683748
$ws.Range("M12").Value = 17.73858496748267
$ws.Range("O12").Value = 2.343858978400021
$ws.Range("C13").Value = 0.1260452914076851
$ws.Range("D13").Value = 0.1122039148211726
$ws.Range("E13").Value = 0.1473537602851636
$ws.Range("F13").Value = 0.9172637347953341
$ws.Range("G13").Value = 0.5278477228691258
$ws.Range("H13").Value = 0.6612966639953584
$ws.Range("I13").Value = 0.6374357782161937
$ws.Range("J13").Value = 0.2624096170303574
$ws.Range("L13").Value = 0.3067082364977693
$ws.Range("M13").Value = 17.64978308985417
$ws.Range("O13").Value = 2.346970093799086
$ws.Range("C14").Value = 0.1268504987179782
$ws.Range("D14").Value = 0.1127536065365433
$ws.Range("E14").Value = 0.1468202008749486
$ws.Range("F14").Value = 0.9223164817277478
$ws.Range("G14").Value = 0.5301686321040506
$ws.Range("H14").Value = 0.6638776031559672
$ws.Range("I14").Value = 0.6408960716931027
$ws.Range("J14").Value = 0.2600126801648202
$ws.Range("L14").Value = 0.3039674639455825
$ws.Range("M14").Value = 17.36012125290398
$ws.Range("O14").Value = 2.357209088815864
$ws.Range("C15").Value = 0.1273471686258603
$ws.Range("D15").Value = 0.1130927028097659
$ws.Range("E15").Value = 0.1464968344996436
$ws.Range("F15").Value = 0.9254398031082758
$ws.Range("G15").Value = 0.5316104667894734
$ws.Range("H15").Value = 0.6654721944853748
$ws.Range("I15").Value = 0.6430356344726604
$ws.Range("J15").Value = 0.2585494847836856
$ws.Range("L15").Value = 0.3022938493841991
$ws.Range("M15").Value = 17.18264842565225
$ws.Range("O15").Value = 2.363552023817931
$ws.Range("C16").Value = 0.1302460193165231
$ws.Range("D16").Value = 0.1150724201841626
$ws.Range("E16").Value = 0.1446960651059825
$ws.Range("F16").Value = 0.9437700838129146
$ws.Range("G16").Value = 0.5401802608582358
$ws.Range("H16").Value = 0.6748189415011296
$ws.Range("I16").Value = 0.6556014964312915
$ws.Range("J16").Value = 0.2502398646578712
$ws.Range("L16").Value = 0.2927813283141631
$ws.Range("M16").Value = 16.16487806119903
$ws.Range("O16").Value = 2.400985311250224
$ws.Range("C17").Value = 0.1320709583660573
$ws.Range("D17").Value = 0.1163192580148333
$ws.Range("E17").Value = 0.1436373581842787
$ws.Range("F17").Value = 0.9553971476961181
$ws.Range("G17").Value = 0.5457089299962306
$ws.Range("H17").Value = 0.6807380199241635
$ws.Range("I17").Value = 0.6635799697197591
$ws.Range("J17").Value = 0.2452081089695071
$ws.Range("L17").Value = 0.2870142583223583
$ws.Range("M17").Value = 15.5398112243534
$ws.Range("O17").Value = 2.424908851718698
$ws.Range("C18").Value = 0.1331375866787639
$ws.Range("D18").Value = 0.1170482072277963
$ws.Range("E18").Value = 0.1430450638561851
$ws.Range("F18").Value = 0.9622238353492421
$ws.Range("G18").Value = 0.5489875155006061
$ws.Range("H18").Value = 0.6842101042648281
$ws.Range("I18").Value = 0.6682671772670545
$ws.Range("J18").Value = 0.2423376930414776
$ws.Range("L18").Value = 0.2837218810151825
$ws.Range("M18").Value = 15.18002056120622
$ws.Range("O18").Value = 2.439018360579738
$ws.Range("C19").Value = 0.1335016295540115
$ws.Range("D19").Value = 0.1172970356999485
$ws.Range("E19").Value = 0.1428473577959224
$ws.Range("F19").Value = 0.9645590049016093
$ws.Range("G19").Value = 0.5501144271138543
$ws.Range("H19").Value = 0.685397257406386
$ws.Range("I19").Value = 0.669870964214148
$ws.Range("J19").Value = 0.2413698455136881
$ws.Range("L19").Value = 0.2826113322513208
$ws.Range("M19").Value = 15.05815434819976
$ws.Range("O19").Value = 2.44385528557045
$ws.Range("C20").Value = 0.1318749306444573
$ws.Range("D20").Value = 0.1161853067163392
$ws.Range("E20").Value = 0.1437483307470728
$ws.Range("F20").Value = 0.9541450008346786
$ws.Range("G20").Value = 0.5451101602666171
$ws.Range("H20").Value = 0.6801009186807789
$ws.Range("I20").Value = 0.6627204635039341
$ws.Range("J20").Value = 0.2457412820134408
$ws.Range("L20").Value = 0.287625606398251
$ws.Range("M20").Value = 15.60637820158223
$ws.Range("O20").Value = 2.422325923193199
$ws.Range("C21").Value = 0.1266132684887253
$ws.Range("D21").Value = 0.1125916493101187
$ws.Range("E21").Value = 0.1469761988331513
$ws.Range("F21").Value = 0.9208264483598541
$ws.Range("G21").Value = 0.5294827067591683
$ws.Range("H21").Value = 0.6631166628885978
$ws.Range("I21").Value = 0.6398755195439421
$ws.Range("J21").Value = 0.2607156741087238
$ws.Range("L21").Value = 0.3047714111748974
$ws.Range("M21").Value = 17.44521158828627
$ws.Range("O21").Value = 2.354186779928312
$ws.Range("C22").Value = 0.1233287995434367
$ws.Range("D22").Value = 0.1103499142782525
$ws.Range("E22").Value = 0.1492396861261369
$ws.Range("F22").Value = 0.9003170610235287
$ws.Range("G22").Value = 0.5201709096336629
$ws.Range("H22").Value = 0.6526280437266223
$ws.Range("I22").Value = 0.625839175761687
$ws.Range("J22").Value = 0.2707245947497796
$ws.Range("L22").Value = 0.3162080012355375
$ws.Range("M22").Value = 18.64487616100968
$ws.Range("O22").Value = 2.312834300966813
$ws.Range("C23").Value = 0.1250675037139786
$ws.Range("D23").Value = 0.1115364958332421
$ws.Range("E23").Value = 0.1480172568426994
$ws.Range("F23").Value = 0.9111461072986486
$ws.Range("G23").Value = 0.5250571544246725
$ws.Range("H23").Value = 0.6581695843393334
$ws.Range("I23").Value = 0.6332478611553682
$ws.Range("J23").Value = 0.2653617493511859
$ws.Range("L23").Value = 0.3100824068281725
$ws.Range("M23").Value = 18.0047723226387
$ws.Range("O23").Value = 2.334610488068336
$ws.Range("C24").Value = 0.1319635004700217
$ws.Range("D24").Value = 0.116245828349804
$ws.Range("E24").Value = 0.1436981091230507
$ws.Range("F24").Value = 0.95471065401928
$ws.Range("G24").Value = 0.5453805525914248
$ws.Range("H24").Value = 0.6803887369439394
$ws.Range("I24").Value = 0.6631087341411117
$ws.Range("J24").Value = 0.2455001648064155
$ws.Range("L24").Value = 0.2873491438261198
$ws.Range("M24").Value = 15.57628461783537
$ws.Range("O24").Value = 2.423492558501039
$ws.Range("C25").Value = 0.1400490336070064
$ws.Range("D25").Value = 0.121775775007805
$ws.Range("E25").Value = 0.1396726351035156
$ws.Range("F25").Value = 1.007004697398941
$ws.Range("G25").Value = 0.571060510862111
$ws.Range("H25").Value = 0.7069340799855084
$ws.Range("I25").Value = 0.6990616424386324
$ws.Range("J25").Value = 0.2249772118407805
$ws.Range("L25").Value = 0.2637659000964305
$ws.Range("M25").Value = 12.94636565306166
$ws.Range("O25").Value = 2.532681409165036
